# Add new columns I (I0) and J (IF) to Sheet1, mirroring columns H values
# (I = min-like helper column, J = H value) for rows 2-49, per commit
# "I0 and IF added".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Apply the same style used by the other header cells (B1:H1) to the new
# header cells so they match (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data rows 2-49 for columns I (I0) and J (IF)
$iValues = @(1,1,1,1,1,1,1,1,4,4,1,1,5,2,7,8,8,9,8,1,5,5,7,9,6,7,8,5,10,7,7,6,8,5,8,7,7,8,5,9,6,1,1,1,1,1,1,3)
$jValues = @(6,6,4,6,7,7,6,5,7,7,4,3,7,3,8,8,9,10,8,1,6,6,7,9,7,8,8,7,10,8,8,7,9,6,8,8,7,8,7,9,7,3,4,4,3,3,2,3)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
